$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Update Day 3 estimates (column H) for task row 3 and the "Skattat" summary row 6
$ws.Range("H3").Value = 3
$ws.Range("H6").Value = 3

# Add a new styled (wrap text, white fill) blank cell below the chart
$ws.Range("F13").Interior.ThemeColor = 2

# Update the selection to mimic the final cursor position left by the author
$ws.Range("H9").Select()
